$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns: "prot" -> "prot_id", "protNew" -> "new_prot_id"
$ws.Range("A1").Value = "prot_id"
$ws.Range("B1").Value = "new_prot_id"

# Update the active selection to B2
$ws.Range("B2").Select()
